# Legs Update and Sesi 1 Update!
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Leg measurement updates (row 6 = L3, row 8 = R2, row 9 = R3) ---
# Dependent formula cells (H, I, K, L columns) recalc automatically.
$ws.Range("G6").Value = 2120

$ws.Range("D8").Value = 1270
$ws.Range("G8").Value = 750

$ws.Range("C9").Value = 1400
$ws.Range("D9").Value = 1250

# --- Sesi 1 Update: move current selection to C14 ---
$ws.Range("C14").Select()
